# Refresh the cryptocurrency snapshot (prices + 1h volume change) pulled in by the
# GitHub Actions scraper job. Only the D (Price) and E (Volume(1h)) columns move;
# coin name/link/rank columns are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Some prices (e.g. "230.37") parse as plain numbers, and a bare
    # Range.Value assignment would let Excel auto-convert them to numeric
    # cells. Forcing Text format for the assignment keeps them as strings
    # (matching the original inline-string cells), then restore the default
    # "Normal" style so no stray per-cell formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "38.249.96"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "2.064.52"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "230.37"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("E6").Value = "  +1.68%  "
Set-TextValue $ws.Range("D7") "58.18"
$ws.Range("E7").Value = "  +6.67%  "
Set-TextValue $ws.Range("D8") "0.998"
$ws.Range("E8").Value = "  -0.30%  "
Set-TextValue $ws.Range("D9") "0.388"
$ws.Range("E9").Value = "  +2.85%  "
Set-TextValue $ws.Range("D10") "0.0807"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "2.364.54"
$ws.Range("E12").Value = "  +2.47%  "
Set-TextValue $ws.Range("D13") "14.62"
$ws.Range("E13").Value = "  +3.83%  "
Set-TextValue $ws.Range("D14") "20.74"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("E15").Value = "  +2.40%  "
Set-TextValue $ws.Range("D16") "5.31"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "2.057.44"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "38.027.19"
Set-TextValue $ws.Range("D19") "6.19"
$ws.Range("E19").Value = "  +1.78%  "
Set-TextValue $ws.Range("D20") "69.83"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("E21").Value = "  +2.03%  "
Set-TextValue $ws.Range("D22") "225.18"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +3.35%  "
Set-TextValue $ws.Range("D26") "9.35"
$ws.Range("E26").Value = "  +2.46%  "
Set-TextValue $ws.Range("D27") "166.17"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +7.84%  "
Set-TextValue $ws.Range("D29") "19.05"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("E31").Value = "  +2.06%  "
Set-TextValue $ws.Range("D32") "4.57"
$ws.Range("E32").Value = "  +1.89%  "
Set-TextValue $ws.Range("D33") "4.62"
$ws.Range("E33").Value = "  +5.07%  "
Set-TextValue $ws.Range("D34") "0.0617"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +6.94%  "
Set-TextValue $ws.Range("D36") "2.39"
$ws.Range("E36").Value = "  +2.56%  "
Set-TextValue $ws.Range("D37") "6.11"
$ws.Range("E37").Value = "  +14.96%  "
$ws.Range("E38").Value = "  +6.25%  "
$ws.Range("E39").Value = "  -0.08%  "
Set-TextValue $ws.Range("D40") "98.55"
$ws.Range("E40").Value = "  +4.33%  "
Set-TextValue $ws.Range("D41") "0.0219"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").Value = "1.486.08"
$ws.Range("E42").Value = "  +0.59%  "
Set-TextValue $ws.Range("D43") "16.96"
$ws.Range("E43").Value = "  +3.24%  "
Set-TextValue $ws.Range("D44") "0.0946"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("E45").Value = "  +4.02%  "
$ws.Range("E46").Value = "  +0.39%  "
Set-TextValue $ws.Range("D47") "4.12"
$ws.Range("E47").Value = "  +18.13%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "2.250.90"
$ws.Range("E51").Value = "  +2.59%  "
